# ---------------------------------------------------------------------------
# Edit summary (see commit diff):
#   1. Slide 16's table (graphic frame "Google Shape;213;p29") switches its
#      table style from the custom "{C6F145C7-89E0-48C5-9310-1F61F265DC4B}"
#      style to the built-in "{339A1028-9747-427F-9DD9-638AEE6A834F}" style.
#   2. The presentation's applied design theme changes from the "Integral"
#      theme colours to the default "Office Theme" colours (dk2/lt2/accent1-6/
#      hlink/folHlink all change; dk1 and lt1 stay black/white in both).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on Slide 16 -------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{339A1028-9747-427F-9DD9-638AEE6A834F}")
    }
}

# --- 2. Swap the deck's theme colours: Integral -> Office Theme -----------
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477    # folHlink -> 954F72
